$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 4).Formula = '="29.709.34"'
$ws.Cells.Item(2, 4).Copy()
$ws.Cells.Item(2, 4).PasteSpecial(-4163)
$ws.Cells.Item(2, 5).Value = '  -0.10%  '
$ws.Cells.Item(3, 4).Formula = '="1.602.95"'
$ws.Cells.Item(3, 4).Copy()
$ws.Cells.Item(3, 4).PasteSpecial(-4163)
$ws.Cells.Item(3, 5).Value = '  -0.77%  '
$ws.Cells.Item(4, 5).Value = '  +0.21%  '
$ws.Cells.Item(5, 4).Formula = '="212.58"'
$ws.Cells.Item(5, 4).Copy()
$ws.Cells.Item(5, 4).PasteSpecial(-4163)
$ws.Cells.Item(5, 5).Value = '  -0.96%  '
$ws.Cells.Item(6, 5).Value = '  +0.29%  '
$ws.Cells.Item(7, 5).Value = '  +0.27%  '
$ws.Cells.Item(8, 4).Formula = '="28.17"'
$ws.Cells.Item(8, 4).Copy()
$ws.Cells.Item(8, 4).PasteSpecial(-4163)
$ws.Cells.Item(8, 5).Value = '  +4.58%  '
$ws.Cells.Item(9, 5).Value = '  +1.20%  '
$ws.Cells.Item(10, 4).Formula = '="0.0603"'
$ws.Cells.Item(10, 4).Copy()
$ws.Cells.Item(10, 4).PasteSpecial(-4163)
$ws.Cells.Item(10, 5).Value = '  +0.63%  '
$ws.Cells.Item(11, 4).Formula = '="0.0908"'
$ws.Cells.Item(11, 4).Copy()
$ws.Cells.Item(11, 4).PasteSpecial(-4163)
$ws.Cells.Item(11, 5).Value = '  -0.76%  '
$ws.Cells.Item(12, 4).Formula = '="1.832.87"'
$ws.Cells.Item(12, 4).Copy()
$ws.Cells.Item(12, 4).PasteSpecial(-4163)
$ws.Cells.Item(12, 5).Value = '  -0.91%  '
$ws.Cells.Item(13, 4).Formula = '="1.600.35"'
$ws.Cells.Item(13, 4).Copy()
$ws.Cells.Item(13, 4).PasteSpecial(-4163)
$ws.Cells.Item(13, 5).Value = '  -0.81%  '
$ws.Cells.Item(14, 4).Formula = '="0.549"'
$ws.Cells.Item(14, 4).Copy()
$ws.Cells.Item(14, 4).PasteSpecial(-4163)
$ws.Cells.Item(14, 5).Value = '  +1.79%  '
$ws.Cells.Item(15, 4).Formula = '="29.740.36"'
$ws.Cells.Item(15, 4).Copy()
$ws.Cells.Item(15, 4).PasteSpecial(-4163)
$ws.Cells.Item(15, 5).Value = '  -0.21%  '
$ws.Cells.Item(16, 5).Value = '  -0.12%  '
$ws.Cells.Item(17, 4).Formula = '="64.13"'
$ws.Cells.Item(17, 4).Copy()
$ws.Cells.Item(17, 4).PasteSpecial(-4163)
$ws.Cells.Item(17, 5).Value = '  +0.71%  '
$ws.Cells.Item(18, 4).Formula = '="241.82"'
$ws.Cells.Item(18, 4).Copy()
$ws.Cells.Item(18, 4).PasteSpecial(-4163)
$ws.Cells.Item(18, 5).Value = '  -1.55%  '
$ws.Cells.Item(19, 5).Value = '  +3.01%  '
$ws.Cells.Item(20, 4).Formula = '="0.0₃0697"'
$ws.Cells.Item(20, 4).Copy()
$ws.Cells.Item(20, 4).PasteSpecial(-4163)
$ws.Cells.Item(20, 5).Value = '  +0.10%  '
$ws.Cells.Item(21, 4).Formula = '="0.999"'
$ws.Cells.Item(21, 4).Copy()
$ws.Cells.Item(21, 4).PasteSpecial(-4163)
$ws.Cells.Item(21, 5).Value = '  +0.58%  '
$ws.Cells.Item(22, 5).Value = '  -1.00%  '
$ws.Cells.Item(23, 4).Formula = '="9.43"'
$ws.Cells.Item(23, 4).Copy()
$ws.Cells.Item(23, 4).PasteSpecial(-4163)
$ws.Cells.Item(23, 5).Value = '  +1.51%  '
$ws.Cells.Item(24, 4).Formula = '="2.11"'
$ws.Cells.Item(24, 4).Copy()
$ws.Cells.Item(24, 4).PasteSpecial(-4163)
$ws.Cells.Item(24, 5).Value = '  -0.42%  '
$ws.Cells.Item(25, 4).Formula = '="154.98"'
$ws.Cells.Item(25, 4).Copy()
$ws.Cells.Item(25, 4).PasteSpecial(-4163)
$ws.Cells.Item(25, 5).Value = '  -0.66%  '
$ws.Cells.Item(26, 5).Value = '  +0.35%  '
$ws.Cells.Item(27, 4).Formula = '="0.109"'
$ws.Cells.Item(27, 4).Copy()
$ws.Cells.Item(27, 4).PasteSpecial(-4163)
$ws.Cells.Item(27, 5).Value = '  +0.39%  '
$ws.Cells.Item(28, 4).Formula = '="6.44"'
$ws.Cells.Item(28, 4).Copy()
$ws.Cells.Item(28, 4).PasteSpecial(-4163)
$ws.Cells.Item(28, 5).Value = '  +0.35%  '
$ws.Cells.Item(29, 5).Value = '  +0.13%  '
$ws.Cells.Item(30, 5).Value = '  +0.87%  '
$ws.Cells.Item(31, 5).Value = '  -0.13%  '
$ws.Cells.Item(32, 5).Value = '  -0.57%  '
$ws.Cells.Item(33, 5).Value = '  +2.18%  '
$ws.Cells.Item(34, 4).Formula = '="1.420.55"'
$ws.Cells.Item(34, 4).Copy()
$ws.Cells.Item(34, 4).PasteSpecial(-4163)
$ws.Cells.Item(34, 5).Value = '  -1.77%  '
$ws.Cells.Item(35, 5).Value = '  +2.84%  '
$ws.Cells.Item(36, 4).Formula = '="2.90"'
$ws.Cells.Item(36, 4).Copy()
$ws.Cells.Item(36, 4).PasteSpecial(-4163)
$ws.Cells.Item(36, 5).Value = '  +1.85%  '
$ws.Cells.Item(37, 5).Value = '  -1.99%  '
$ws.Cells.Item(38, 5).Value = '  -0.95%  '
$ws.Cells.Item(39, 5).Value = '  +1.08%  '
$ws.Cells.Item(40, 4).Formula = '="0.546"'
$ws.Cells.Item(40, 4).Copy()
$ws.Cells.Item(40, 4).PasteSpecial(-4163)
$ws.Cells.Item(40, 5).Value = '  +1.28%  '
$ws.Cells.Item(41, 4).Formula = '="56.19"'
$ws.Cells.Item(41, 4).Copy()
$ws.Cells.Item(41, 4).PasteSpecial(-4163)
$ws.Cells.Item(41, 5).Value = '  -0.42%  '
$ws.Cells.Item(42, 5).Value = '  +5.66%  '
$ws.Cells.Item(43, 4).Formula = '="0.816"'
$ws.Cells.Item(43, 4).Copy()
$ws.Cells.Item(43, 4).PasteSpecial(-4163)
$ws.Cells.Item(43, 5).Value = '  +1.83%  '
$ws.Cells.Item(44, 5).Value = '  -0.62%  '
$ws.Cells.Item(45, 5).Value = '  +0.19%  '
$ws.Cells.Item(46, 4).Formula = '="67.34"'
$ws.Cells.Item(46, 4).Copy()
$ws.Cells.Item(46, 4).PasteSpecial(-4163)
$ws.Cells.Item(46, 5).Value = '  -2.95%  '
$ws.Cells.Item(47, 4).Formula = '="0.985"'
$ws.Cells.Item(47, 4).Copy()
$ws.Cells.Item(47, 4).PasteSpecial(-4163)
$ws.Cells.Item(47, 5).Value = '  +17.52%  '
$ws.Cells.Item(48, 4).Formula = '="5.37"'
$ws.Cells.Item(48, 4).Copy()
$ws.Cells.Item(48, 4).PasteSpecial(-4163)
$ws.Cells.Item(48, 5).Value = '  +0.31%  '
$ws.Cells.Item(49, 4).Formula = '="1.741.36"'
$ws.Cells.Item(49, 4).Copy()
$ws.Cells.Item(49, 4).PasteSpecial(-4163)
$ws.Cells.Item(49, 5).Value = '  -0.96%  '
$ws.Cells.Item(50, 4).Formula = '="86.59"'
$ws.Cells.Item(50, 4).Copy()
$ws.Cells.Item(50, 4).PasteSpecial(-4163)
$ws.Cells.Item(51, 2).Value = 'Cronos'
$ws.Cells.Item(51, 3).Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Cells.Item(51, 4).Formula = '="0.0524"'
$ws.Cells.Item(51, 4).Copy()
$ws.Cells.Item(51, 4).PasteSpecial(-4163)
$ws.Cells.Item(51, 5).Value = '  +0.57%  '
$excel.CutCopyMode = $false
